# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'28.203.59"
$ws.Cells.Item(2, 5).Value = '  -5.39%  '

$ws.Cells.Item(3, 4).Value = "'1.836.69"
$ws.Cells.Item(3, 5).Value = '  -5.04%  '

$ws.Cells.Item(4, 5).Value = '  -0.64%  '

$ws.Cells.Item(5, 4).Value = "'330.23"
$ws.Cells.Item(5, 5).Value = '  -2.21%  '

$ws.Cells.Item(6, 5).Value = '  -0.54%  '

$ws.Cells.Item(7, 4).Value = "'0.4606"
$ws.Cells.Item(7, 5).Value = '  -4.69%  '

$ws.Cells.Item(8, 4).Value = "'0.3857"
$ws.Cells.Item(8, 5).Value = '  -6.12%  '

$ws.Cells.Item(9, 4).Value = "'46.03"
$ws.Cells.Item(9, 5).Value = '  -3.88%  '

$ws.Cells.Item(10, 4).Value = "'0.07850"
$ws.Cells.Item(10, 5).Value = '  -3.82%  '

$ws.Cells.Item(11, 5).Value = '  -5.60%  '

$ws.Cells.Item(12, 4).Value = "'21.90"
$ws.Cells.Item(12, 5).Value = '  -7.43%  '

$ws.Cells.Item(13, 4).Value = "'1.831.60"
$ws.Cells.Item(13, 5).Value = '  -7.99%  '

$ws.Cells.Item(14, 4).Value = "'5.701"
$ws.Cells.Item(14, 5).Value = '  -6.12%  '

$ws.Cells.Item(15, 5).Value = '  -5.08%  '

$ws.Cells.Item(16, 4).Value = "'0.06850"
$ws.Cells.Item(16, 5).Value = '  +0.17%  '

$ws.Cells.Item(17, 4).Value = "'1.002"
$ws.Cells.Item(17, 5).Value = '  -0.66%  '

$ws.Cells.Item(18, 4).Value = "'86.81"
$ws.Cells.Item(18, 5).Value = '  -4.41%  '

$ws.Cells.Item(19, 4).Value = "'0.000009941"
$ws.Cells.Item(19, 5).Value = '  -3.81%  '

$ws.Cells.Item(20, 4).Value = "'16.91"
$ws.Cells.Item(20, 5).Value = '  -4.76%  '

$ws.Cells.Item(21, 4).Value = "'1.001"
$ws.Cells.Item(21, 5).Value = '  -0.57%  '

$ws.Cells.Item(22, 4).Value = "'28.230.42"
$ws.Cells.Item(22, 5).Value = '  -5.30%  '

$ws.Cells.Item(23, 5).Value = '  -5.25%  '

$ws.Cells.Item(24, 4).Value = "'10.94"
$ws.Cells.Item(24, 5).Value = '  -7.94%  '

$ws.Cells.Item(25, 4).Value = "'2.131"
$ws.Cells.Item(25, 5).Value = '  -2.20%  '

$ws.Cells.Item(26, 4).Value = "'2.043.22"
$ws.Cells.Item(26, 5).Value = '  -3.57%  '

$ws.Cells.Item(27, 4).Value = "'153.27"
$ws.Cells.Item(27, 5).Value = '  -2.37%  '

$ws.Cells.Item(28, 4).Value = "'19.19"

$ws.Cells.Item(29, 4).Value = "'5.691"
$ws.Cells.Item(29, 5).Value = '  -13.35%  '

$ws.Cells.Item(30, 4).Value = "'1.974"
$ws.Cells.Item(30, 5).Value = '  -5.45%  '

$ws.Cells.Item(31, 4).Value = "'116.61"
$ws.Cells.Item(31, 5).Value = '  -3.55%  '

$ws.Cells.Item(32, 2).Value = 'Stellar'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(32, 4).Value = "'0.09278"
$ws.Cells.Item(32, 5).Value = '  -3.90%  '

$ws.Cells.Item(33, 2).Value = 'ImmutableX'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(33, 4).Value = "'0.9371"
$ws.Cells.Item(33, 5).Value = '  -6.91%  '

$ws.Cells.Item(34, 4).Value = "'5.270"
$ws.Cells.Item(34, 5).Value = '  -4.85%  '

$ws.Cells.Item(35, 4).Value = "'3.447"
$ws.Cells.Item(35, 5).Value = '  -2.40%  '

$ws.Cells.Item(36, 5).Value = '  -6.27%  '

$ws.Cells.Item(37, 4).Value = "'0.05997"
$ws.Cells.Item(37, 5).Value = '  -8.88%  '

$ws.Cells.Item(38, 4).Value = "'0.02147"
$ws.Cells.Item(38, 5).Value = '  -5.87%  '

$ws.Cells.Item(39, 4).Value = "'1.146"
$ws.Cells.Item(39, 5).Value = '  -4.53%  '

$ws.Cells.Item(40, 5).Value = '  -0.56%  '

$ws.Cells.Item(41, 4).Value = "'7.602"
$ws.Cells.Item(41, 5).Value = '  -4.16%  '

$ws.Cells.Item(42, 4).Value = "'0.5610"
$ws.Cells.Item(42, 5).Value = '  -5.95%  '

$ws.Cells.Item(43, 4).Value = "'9.972"
$ws.Cells.Item(43, 5).Value = '  -7.29%  '

$ws.Cells.Item(44, 4).Value = "'0.1771"
$ws.Cells.Item(44, 5).Value = '  -4.05%  '

$ws.Cells.Item(45, 4).Value = "'1.250"
$ws.Cells.Item(45, 5).Value = '  -1.79%  '

$ws.Cells.Item(46, 4).Value = "'2.243"
$ws.Cells.Item(46, 5).Value = '  -9.15%  '

$ws.Cells.Item(47, 4).Value = "'11.56"
$ws.Cells.Item(47, 5).Value = '  -5.82%  '

$ws.Cells.Item(48, 4).Value = "'0.5269"
$ws.Cells.Item(48, 5).Value = '  -5.08%  '

$ws.Cells.Item(49, 4).Value = "'0.07015"
$ws.Cells.Item(49, 5).Value = '  -6.24%  '

$ws.Cells.Item(50, 4).Value = "'1.829"
$ws.Cells.Item(50, 5).Value = '  -7.65%  '

$ws.Cells.Item(51, 4).Value = "'112.54"
$ws.Cells.Item(51, 5).Value = '  -3.64%  '
